$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Programs")

# ---------------------------------------------------------------------------
# 1) Establish the new cell style (numFmt=General, fontId=5, left/vcenter)
#    on E124 -- a cell that legitimately needs it for the new program below.
#    Seed it from an existing "comment" style (fontId 5, used on column E)
#    and then tweak the alignment so the engine allocates exactly one new
#    cellXfs entry.
# ---------------------------------------------------------------------------
$ws.Range("E85").Copy()
$ws.Range("E124").PasteSpecial(-4122)            # xlPasteFormats
$ws.Range("E124").HorizontalAlignment = -4131    # xlLeft
$ws.Range("E124").VerticalAlignment = -4108      # xlCenter

# Propagate that same new style to every other cell that needs it.
$ws.Range("E124").Copy()
$ws.Range("E88").PasteSpecial(-4122)
$ws.Range("E89").PasteSpecial(-4122)
$ws.Range("E90").PasteSpecial(-4122)
$ws.Range("E125").PasteSpecial(-4122)
$ws.Range("E126").PasteSpecial(-4122)
$ws.Range("E129").PasteSpecial(-4122)
$ws.Range("E130").PasteSpecial(-4122)
$ws.Range("E131").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the existing "Increments A/B" comment block (rows 88-90): the
#    comment now belongs to the merged E88:E90 block instead of just E89.
# ---------------------------------------------------------------------------
$ws.Range("E88").Value = "// Increments B"
$ws.Range("E89").ClearContents()
$ws.Range("E90").ClearContents()
$ws.Range("E88:E90").Merge()

# ---------------------------------------------------------------------------
# 3) New "copy a block of memory" test program, rows 118-132.
# ---------------------------------------------------------------------------

# -- Column A (source addresses used by the DEC2HEX formula in column B) --
$aValues = @{
    118 = 0;  119 = 3;  120 = 6;  121 = 9;  122 = 12; 123 = 15; 124 = 18;
    125 = 21; 126 = 24; 127 = 27; 128 = 30; 129 = 33; 130 = 36; 131 = 39; 132 = 42
}
foreach ($r in $aValues.Keys) {
    $ws.Cells.Item($r, 1).Value = $aValues[$r]
}

# -- Column B: shared "0x" & DEC2HEX(A,3) formulas, same grouping as the diff --
$ws.Range("B118:B123").Formula = '="0x" & DEC2HEX(A118,3)'
$ws.Range("B123:B132").Formula = '="0x" & DEC2HEX(A123,3)'

# -- Column C (mnemonics) --
$ws.Range("C118").Value = "LD B, 0x30"
$ws.Range("C119").Value = "LD F, 0x3c"
$ws.Range("C120").Value = "LD C, 0x60"
$ws.Range("C121").Value = "LD E, 0x1"
$ws.Range("C122").Value = "LD A, [?B]"
$ws.Range("C123").Value = "ST [?C], A"
$ws.Range("C124").Value = "LD A, B"
$ws.Range("C125").Value = "ADD A, E"
$ws.Range("C126").Value = "LD B, A"
$ws.Range("C127").Value = "SUB A, F"
$ws.Range("C128").Value = "JP Z, [0x01e]"
$ws.Range("C129").Value = "LD A, C"
$ws.Range("C130").Value = "ADD A, E"
$ws.Range("C131").Value = "LD C, A"
$ws.Range("C132").Value = "JP [0x00c]"

# -- Column D (raw bytes; not every row has one) --
$ws.Range("D118").Value = "04 80 30"
$ws.Range("D119").Value = "06 80 3c"
$ws.Range("D121").Value = "06 00 01"
$ws.Range("D122").Value = "10 10 00"
$ws.Range("D124").Value = "08 10 00"
$ws.Range("D125").Value = "14 40 00"
$ws.Range("D126").Value = "08 80 00"
$ws.Range("D127").Value = "18 50 00"
$ws.Range("D128").Value = "30 00 1e"
$ws.Range("D130").Value = "14 40 00"
$ws.Range("D132").Value = "2c 00 0c"

# -- Column E (comments); cells that are part of a merge stay blank --
$ws.Range("E118").Value = "// Source addr of first byte"
$ws.Range("E119").Value = "// Source addr of last byte + 1"
$ws.Range("E120").Value = "// Destitation addr of first byte"
$ws.Range("E121").Value = "// Constant to use in increment"
$ws.Range("E122").Value = "// Reads byte from memory"
$ws.Range("E123").Value = "// Store byte in destiny addr"
$ws.Range("E124").Value = "// Increments B"
$ws.Range("E127").Value = "// Tests whether is the last char"
$ws.Range("E128").Value = "// If true, stops here"
$ws.Range("E129").Value = "// Increments C"
$ws.Range("E132").Value = "// If false, next loop iteration"

# -- style E118..E123, E127, E128, E132 like the rest of the "comment" column --
$ws.Range("E85").Copy()
$ws.Range("E118").PasteSpecial(-4122)
$ws.Range("E119").PasteSpecial(-4122)
$ws.Range("E120").PasteSpecial(-4122)
$ws.Range("E121").PasteSpecial(-4122)
$ws.Range("E122").PasteSpecial(-4122)
$ws.Range("E123").PasteSpecial(-4122)
$ws.Range("E127").PasteSpecial(-4122)
$ws.Range("E128").PasteSpecial(-4122)
$ws.Range("E132").PasteSpecial(-4122)

# -- merge the new comment blocks --
$ws.Range("E124:E126").Merge()
$ws.Range("E129:E131").Merge()

# ---------------------------------------------------------------------------
# 4) View state: selection on the newly-added merged block, scrolled so the
#    new program is visible.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E129:E131").Select()
$excel.ActiveWindow.ScrollRow = 85
